$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.927.83"
$ws.Range("E2").Value = "  -0.35%  "
$ws.Range("D3").Value = "1.910.68"
$ws.Range("E3").Value = "  -0.08%  "
$ws.Range("D4").Value = "'0.9978"
$ws.Range("D5").Value = "'313.15"
$ws.Range("E5").Value = "  -0.78%  "
$ws.Range("D6").Value = "'0.9980"
$ws.Range("E6").Value = "  -0.73%  "
$ws.Range("D7").Value = "'0.5014"
$ws.Range("E7").Value = "  +3.98%  "
$ws.Range("D8").Value = "'0.3818"
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").Value = "'0.07315"
$ws.Range("E9").Value = "  -0.61%  "
$ws.Range("D10").Value = "'0.9120"
$ws.Range("E10").Value = "  -2.38%  "
$ws.Range("D11").Value = "'21.25"
$ws.Range("E11").Value = "  +2.19%  "
$ws.Range("D12").Value = "'0.07669"
$ws.Range("E12").Value = "  -1.83%  "
$ws.Range("D13").Value = "1.913.42"
$ws.Range("E13").Value = "  +1.01%  "
$ws.Range("D14").Value = "'5.483"
$ws.Range("E14").Value = "  -0.42%  "
$ws.Range("D15").Value = "'92.87"
$ws.Range("E15").Value = "  +1.03%  "
$ws.Range("D16").Value = "'0.9986"
$ws.Range("E16").Value = "  -0.77%  "
$ws.Range("D17").Value = "'0.000008744"
$ws.Range("E17").Value = "  -1.52%  "
$ws.Range("D18").Value = "'0.9982"
$ws.Range("E18").Value = "  -0.68%  "
$ws.Range("D19").Value = "27.959.56"
$ws.Range("E19").Value = "  -0.35%  "
$ws.Range("D20").Value = "'14.70"
$ws.Range("D21").Value = "'5.185"
$ws.Range("D22").Value = "2.164.21"
$ws.Range("E22").Value = "  +1.08%  "
$ws.Range("D23").Value = "'10.86"
$ws.Range("E23").Value = "  -0.44%  "
$ws.Range("D24").Value = "'6.617"
$ws.Range("E24").Value = "  -0.27%  "
$ws.Range("D25").Value = "'153.21"
$ws.Range("E25").Value = "  -2.43%  "
$ws.Range("E26").Value = "  -3.86%  "
$ws.Range("D27").Value = "'2.209"
$ws.Range("E27").Value = "  +3.50%  "
$ws.Range("D28").Value = "'18.42"
$ws.Range("E28").Value = "  -0.49%  "
$ws.Range("D29").Value = "'115.44"
$ws.Range("E29").Value = "  -1.30%  "
$ws.Range("D30").Value = "'4.933"
$ws.Range("E30").Value = "  -0.80%  "
$ws.Range("D31").Value = "'0.09024"
$ws.Range("E31").Value = "  +0.78%  "
$ws.Range("D32").Value = "'3.207"
$ws.Range("E32").Value = "  -2.44%  "
$ws.Range("D33").Value = "'4.858"
$ws.Range("E33").Value = "  +4.05%  "
$ws.Range("E34").Value = "  -1.46%  "
$ws.Range("D35").Value = "'0.7799"
$ws.Range("E35").Value = "  +0.73%  "
$ws.Range("B36").Value = "RenderToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D36").Value = "'2.614"
$ws.Range("E36").Value = "  +0.26%  "
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").Value = "'0.02084"
$ws.Range("E37").Value = "  +1.79%  "
$ws.Range("D38").Value = "'3.066"
$ws.Range("E38").Value = "  +2.31%  "
$ws.Range("E39").Value = "  -1.14%  "
$ws.Range("D40").Value = "'0.5559"
$ws.Range("E40").Value = "  +0.77%  "
$ws.Range("D41").Value = "'0.05288"
$ws.Range("E41").Value = "  -0.29%  "
$ws.Range("D42").Value = "'6.881"
$ws.Range("E42").Value = "  -2.19%  "
$ws.Range("D43").Value = "'113.80"
$ws.Range("E43").Value = "  +4.72%  "
$ws.Range("D44").Value = "'8.538"
$ws.Range("E44").Value = "  +0.41%  "
$ws.Range("D45").Value = "'0.1518"
$ws.Range("E45").Value = "  -0.63%  "
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").Value = "'0.4838"
$ws.Range("E46").Value = "  +0.15%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "'10.57"
$ws.Range("E47").Value = "  -0.89%  "
$ws.Range("D48").Value = "'0.9977"
$ws.Range("E48").Value = "  -0.81%  "
$ws.Range("D49").Value = "'1.643"
$ws.Range("E49").Value = "  -0.41%  "
$ws.Range("D50").Value = "'67.60"
$ws.Range("E50").Value = "  -0.49%  "
$ws.Range("D51").Value = "'0.06048"
$ws.Range("E51").Value = "  -0.50%  "
